# Update the "ランサーズ" sheet: refresh rows 2-10 with newly scraped
# listings (timestamp 2025-09-18 06:25:02) and drop the old rows 11-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- 1. Update rows 2-10 in place -----------------------------------------

$ws.Range("A2").Value = "2025-09-18 06:25:02"
$ws.Range("B2").Value = "AIチャットボットのβ版テスト参加者募集!"
$ws.Range("D2").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5395790"
$ws.Range("G2").Value = 295
$ws.Range("H2").Value = "🔥AI,Ai"

$ws.Range("A3").Value = "2025-09-18 06:25:02"
$ws.Range("B3").Value = "【急募】入力ミス防止のためのkintone Javascript開発者募集"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5395671"
$ws.Range("G3").Value = 128
$ws.Range("H3").Value = "★Java ◆開発"

$ws.Range("A4").Value = "2025-09-18 06:25:02"
$ws.Range("B4").Value = "初回 Laravel Livewireを使ったWebシステム開発の募集"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5395502"
$ws.Range("G4").Value = 125
$ws.Range("H4").Value = "◆開発,システム開発"

$ws.Range("A5").Value = "2025-09-18 06:25:02"
$ws.Range("B5").Value = "2026年度新入社員研修Javaサブ講師"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5395710"
$ws.Range("G5").Value = 85
$ws.Range("H5").Value = "★Java"

$ws.Range("A6").Value = "2025-09-18 06:25:02"
$ws.Range("B6").Value = "2026年度新入社員研修Javaサブ講師 (4~6月)"
$ws.Range("D6").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5395713"
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = "★Java"

$ws.Range("A7").Value = "2025-09-18 06:25:02"
$ws.Range("B7").Value = "【急募】ストアーズ連携LINE予約サイト制作のプロを探しています!"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5395825"
$ws.Range("G7").Value = 33
$ws.Range("H7").Value = "◇サイト"

$ws.Range("A8").Value = "2025-09-18 06:25:02"
$ws.Range("B8").Value = "【継続案件|お気軽にご応募ください】WebシステムのQAエンジニア募集"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5395741"
$ws.Range("G8").Value = 40
$ws.Range("H8").ClearContents()

$ws.Range("A9").Value = "2025-09-18 06:25:02"
$ws.Range("B9").Value = "【急募】JotformとGoogleスプレッドシート連携のエキスパート募集!"
$ws.Range("D9").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5395809"
$ws.Range("G9").Value = 10
$ws.Range("H9").ClearContents()

$ws.Range("A10").Value = "2025-09-18 06:25:02"
$ws.Range("B10").Value = "【MT4】ゴールドの取引を行うEAのサンプルソース納品"
$ws.Range("D10").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5395799"
$ws.Range("G10").Value = 10
$ws.Range("H10").ClearContents()

# --- 2. Drop the now-stale rows 11-23 --------------------------------------

$ws.Rows("11:23").Delete()

# --- 3. Rebuild the hyperlinks for F2:F10 with the refreshed URLs ---------

$ws.Range("F2:F10").Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5395790"
    3  = "https://www.lancers.jp/work/detail/5395671"
    4  = "https://www.lancers.jp/work/detail/5395502"
    5  = "https://www.lancers.jp/work/detail/5395710"
    6  = "https://www.lancers.jp/work/detail/5395713"
    7  = "https://www.lancers.jp/work/detail/5395825"
    8  = "https://www.lancers.jp/work/detail/5395741"
    9  = "https://www.lancers.jp/work/detail/5395809"
    10 = "https://www.lancers.jp/work/detail/5395799"
}

foreach ($r in 2..10) {
    $cell = $ws.Range("F" + $r)
    $ws.Hyperlinks.Add($cell, $urls[$r])
    $cell.Style = "Hyperlink"
}

# --- 4. Resize columns B and H ----------------------------------------------
# ColumnWidth uses "character" units which Excel offsets by 5/6 of a
# character from the raw OOXML <col width> value, so back that offset out
# to land exactly on width=40 / width=12 in the saved file.

$ws.Columns("B").ColumnWidth = 40 - 5/6
$ws.Columns("H").ColumnWidth = 12 - 5/6
